$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1877
$ws.Range("L3").Value = 1895
$ws.Range("C4").Value = 1868
$ws.Range("K4").Value = 1756
$ws.Range("L4").Value = 534
$ws.Range("L5").Value = 113
$ws.Range("L6").Value = 1736
$ws.Range("C7").Value = 28412
$ws.Range("K7").Value = 27547
$ws.Range("L7").Value = 6155

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 106
$ws.Range("L3").Value = 130
$ws.Range("L5").Value = 15
$ws.Range("L7").Value = 384

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 46
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 148

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 65
$ws.Range("L3").Value = 92
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 270

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 72
$ws.Range("L3").Value = 61
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 208
$ws.Range("L8").Value = 384
$ws.Range("L11").Value = 113
$ws.Range("L15").Value = 43
$ws.Range("L19").Value = 175
$ws.Range("L20").Value = 162
$ws.Range("L21").Value = 20
$ws.Range("L23").Value = 67
$ws.Range("L29").Value = 315
$ws.Range("L31").Value = 61
$ws.Range("L33").Value = 270
$ws.Range("L34").Value = 38
$ws.Range("L36").Value = 89
$ws.Range("L37").Value = 222
$ws.Range("L44").Value = 41
$ws.Range("L48").Value = 91
$ws.Range("L49").Value = 32
$ws.Range("L51").Value = 74
$ws.Range("C63").Value = 293
$ws.Range("K63").Value = 87
$ws.Range("L63").Value = 17
$ws.Range("L64").Value = 44
$ws.Range("L67").Value = 215
$ws.Range("L75").Value = 25
$ws.Range("L77").Value = 39
$ws.Range("L79").Value = 169
$ws.Range("L83").Value = 148
$ws.Range("L84").Value = 62
$ws.Range("L85").Value = 325
$ws.Range("L86").Value = 46
$ws.Range("L91").Value = 85
$ws.Range("L94").Value = 75
$ws.Range("L99").Value = 96
$ws.Range("C101").Value = 28412
$ws.Range("K101").Value = 27547
$ws.Range("L101").Value = 6155

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 64
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 215

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 112
$ws.Range("L7").Value = 315

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 56
$ws.Range("L3").Value = 53
$ws.Range("L6").Value = 56
$ws.Range("L7").Value = 175

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 22
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 34
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 61
$ws.Range("L3").Value = 55
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L3").Value = 10
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 162

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 89

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 65
$ws.Range("L4").Value = 20
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 208

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 21
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 17
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 39
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 46

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 25

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 99
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 325

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 39
